$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, pushing existing rows 74..85 down to 75..86.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the weekly record.
$ws.Cells.Item(74, 1).Value = 4
$ws.Cells.Item(74, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(74, 3).Value = "Los Lagos"
$ws.Cells.Item(74, 4).Value = 44551
$ws.Cells.Item(74, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(74, 5).Value = 10
$ws.Cells.Item(74, 6).Value = 100112052
$ws.Cells.Item(74, 7).Value = "Albahaca"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 80
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 8000
$ws.Cells.Item(74, 14).Value = "`$/docena de matas"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 1333
$ws.Cells.Item(74, 17).Value = 6
$ws.Cells.Item(74, 18).Value = "Hortaliza"
